$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (Price / Volume) to Text format before writing,
# so Excel does not auto-coerce numeric-looking strings (e.g. "1.00", "0.617")
# into actual numbers and lose formatting / type fidelity.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
  "D2" = '67.827.42'
  "E2" = '  +0.86%  '
  "D3" = '3.904.86'
  "E3" = '  +2.03%  '
  "E4" = '  -0.07%  '
  "D5" = '478.02'
  "E5" = '  +4.72%  '
  "D6" = '148.09'
  "E6" = '  +0.58%  '
  "D7" = '0.617'
  "E7" = '  -0.88%  '
  "E8" = '  -0.07%  '
  "E9" = '  -2.28%  '
  "E10" = '  +8.32%  '
  "D11" = '0.0000352'
  "E11" = '  +11.32%  '
  "D12" = '42.49'
  "E12" = '  -2.71%  '
  "D13" = '4.513.60'
  "E13" = '  +1.88%  '
  "D14" = '10.26'
  "E14" = '  -0.97%  '
  "D15" = '14.80'
  "E15" = '  -1.14%  '
  "D16" = '3.963.03'
  "E16" = '  +4.35%  '
  "D17" = '0.137'
  "E17" = '  -0.23%  '
  "D18" = '19.88'
  "E18" = '  -1.15%  '
  "D19" = '1.13'
  "E19" = '  -3.10%  '
  "D20" = '67.961.91'
  "E20" = '  +0.96%  '
  "D21" = '433.29'
  "E21" = '  +0.32%  '
  "D22" = '3.34'
  "E22" = '  +2.13%  '
  "D23" = '14.39'
  "E23" = '  -2.43%  '
  "D24" = '87.31'
  "E24" = '  +1.01%  '
  "D25" = '3.54'
  "E25" = '  +1.23%  '
  "B26" = 'RenderToken'
  "C26" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
  "D26" = '10.47'
  "E26" = '  +1.43%  '
  "B27" = 'EthereumClassic'
  "C27" = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
  "D27" = '38.13'
  "E27" = '  +2.45%  '
  "B28" = 'LEO'
  "C28" = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
  "D28" = '5.87'
  "E28" = '  +6.22%  '
  "B29" = 'Filecoin'
  "C29" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
  "D29" = '10.09'
  "E29" = '  +3.80%  '
  "B30" = 'Bittensor'
  "C30" = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
  "D30" = '730.63'
  "E30" = '  -0.52%  '
  "B31" = 'Cosmos'
  "C31" = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
  "D31" = '13.23'
  "E31" = '  -4.26%  '
  "B32" = 'Hedera'
  "C32" = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
  "D32" = '0.128'
  "E32" = '  -5.02%  '
  "B33" = 'Toncoin'
  "C33" = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
  "D33" = '2.81'
  "E33" = '  +2.83%  '
  "D34" = '0.0₃0906'
  "E34" = '  +31.78%  '
  "B35" = 'InjectiveProtocol'
  "C35" = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
  "D35" = '41.95'
  "E35" = '  -2.25%  '
  "B36" = 'OKB'
  "C36" = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
  "D36" = '58.71'
  "E36" = '  +2.27%  '
  "B37" = 'Kaspa'
  "C37" = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
  "D37" = '0.151'
  "E37" = '  -5.03%  '
  "B38" = 'EnergySwap'
  "C38" = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
  "D38" = '33.18'
  "E38" = '  +33.38%  '
  "B39" = 'Dai'
  "C39" = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
  "D39" = '0.999'
  "E39" = '  -0.04%  '
  "B40" = 'NEARProtocol'
  "C40" = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
  "D40" = '5.40'
  "E40" = '  -2.79%  '
  "B41" = 'VeChain'
  "C41" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
  "D41" = '0.0469'
  "E41" = '  -1.27%  '
  "B42" = 'Fetch.AI'
  "C42" = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
  "D42" = '2.82'
  "E42" = '  +7.05%  '
  "D43" = '2.99'
  "E43" = '  +11.79%  '
  "B44" = 'ThetaToken'
  "C44" = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
  "D44" = '2.98'
  "E44" = '  +2.40%  '
  "B45" = 'TheGraph'
  "C45" = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
  "D45" = '0.343'
  "E45" = '  -2.85%  '
  "D46" = '1.00'
  "E46" = '  -0.01%  '
  "B47" = 'Stellar'
  "C47" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
  "D47" = '0.140'
  "E47" = '  +0.10%  '
  "B48" = 'LidoDAOToken'
  "C48" = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
  "D48" = '3.44'
  "E48" = '  -0.58%  '
  "B49" = 'ARBITRUM'
  "C49" = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
  "D49" = '2.17'
  "E49" = '  +1.66%  '
  "D50" = '145.59'
  "E50" = '  +1.37%  '
  "B51" = 'ApeXProtocol'
  "C51" = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
  "D51" = '3.15'
  "E51" = '  -2.47%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# Restore the default (no explicit number format) style on the touched range
# so cells match the original "General" styling once text values are locked in.
$ws.Range("D2:E51").Style = "Normal"
